$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MENTORS")

# Update leadership titles to reflect the 2023 term instead of "Ex-" titles
$ws.Range("B2").Value = "Team Lead 2023"
$ws.Range("B3").Value = "SAMBED Lead 2023"
$ws.Range("B5").Value = "SPACED Lead 2023"
$ws.Range("B4").Value = "SIESED Lead 2023"

# Leave the final selection on C6, matching where the author's cursor ended up
$ws.Range("C6").Select()
